$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(3.3152176530624358, 1.5772195816834524, -0.9064625844065658, 2.769665179102258)"
$ws.Range("C2").Value = "NCT(2.5257719539080146, 1.956392391595798, -0.013672115746600228, 3.9179388709047)"
$ws.Range("D2").Value = "JSU(-1.0371755464286299, 1.2822633889176782, 0.6948798610497393, 3.4111443054072828)"
$ws.Range("E2").Value = "JSU(-1.1384733287182143, 1.045854380514267, 2.8772947593622646, 3.938147739889253)"
